$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value2 = 784.6
$ws.Range("I38").Value2 = 784.6
$ws.Range("K38").Value2 = 2353.8
$ws.Range("M38").Value2 = -1981.8
$ws.Range("H43").Value2 = 2637300.2
$ws.Range("I43").Value2 = 4391047
$ws.Range("J43").Value2 = 6680.375
$ws.Range("K43").Value2 = 4391047
$ws.Range("L43").Value2 = 6680.375
$ws.Range("M43").Value2 = -4390978
$ws.Range("N43").Value2 = -6818.375
$ws.Range("H107").Value2 = 655.871
$ws.Range("I107").Value2 = 674.34784
$ws.Range("J107").Value2 = 602.75
$ws.Range("K107").Value2 = 674.34784
$ws.Range("L107").Value2 = 602.75
$ws.Range("M107").Value2 = 1245.65216
$ws.Range("N107").Value2 = -4442.75
$ws.Range("H113").Value2 = 7863
$ws.Range("I113").Value2 = 3115.6667
$ws.Range("K113").Value2 = 3115.6667
$ws.Range("M113").Value2 = 138.3332999999998
$ws.Range("H121").Value2 = 3356.4285
$ws.Range("I121").Value2 = 0
$ws.Range("K121").Value2 = 0
$ws.Range("M121").ClearContents()
$ws.Range("H132").Value2 = 37041836
$ws.Range("J132").Value2 = 3565.625
$ws.Range("L132").Value2 = 10696.875
$ws.Range("N132").Value2 = -15756.875
$ws.Range("H138").Value2 = 308347.5
$ws.Range("J138").Value2 = 363118.8
$ws.Range("L138").Value2 = 1089356.4
$ws.Range("N138").Value2 = -1099636.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 2416.5334
$ws.Range("I45").Value2 = 1856.5
$ws.Range("J45").Value2 = 6056.75
$ws.Range("K45").Value2 = 1856.5
$ws.Range("L45").Value2 = 6056.75
$ws.Range("M45").Value2 = -1479.5
$ws.Range("N45").Value2 = -6810.75
$ws.Range("H102").Value2 = 2110.2727
$ws.Range("I102").Value2 = 762.75
$ws.Range("K102").Value2 = 762.75
$ws.Range("M102").Value2 = 859.25
$ws.Range("H117").Value2 = 59499.5
$ws.Range("J117").Value2 = 59499.5
$ws.Range("L117").Value2 = 59499.5
$ws.Range("N117").Value2 = -68677.5
$ws.Range("H132").Value2 = 8040.533
$ws.Range("I132").Value2 = 8382.75
$ws.Range("K132").Value2 = 25148.25
$ws.Range("M132").Value2 = -22618.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value2 = 461.53845
$ws.Range("I94").Value2 = 488.63635
$ws.Range("J94").Value2 = 312.5
$ws.Range("K94").Value2 = 488.63635
$ws.Range("L94").Value2 = 312.5
$ws.Range("M94").Value2 = -37.63634999999999
$ws.Range("N94").Value2 = -1214.5
$ws.Range("H117").Value2 = 49000
$ws.Range("J117").Value2 = 49000
$ws.Range("L117").Value2 = 49000
$ws.Range("N117").Value2 = -58178
$ws.Range("H134").Value2 = 4606.5483
$ws.Range("I134").Value2 = 4703.207
$ws.Range("K134").Value2 = 14109.621
$ws.Range("M134").Value2 = -11574.621
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1819.6562
$ws.Range("I31").Value2 = 1261.1428
$ws.Range("J31").Value2 = 2885.9092
$ws.Range("K31").Value2 = 1261.1428
$ws.Range("L31").Value2 = 2885.9092
$ws.Range("M31").Value2 = -966.1428000000001
$ws.Range("N31").Value2 = -3475.9092
$ws.Range("H34").Value2 = 1819.6562
$ws.Range("I34").Value2 = 1261.1428
$ws.Range("J34").Value2 = 2885.9092
$ws.Range("K34").Value2 = 1261.1428
$ws.Range("L34").Value2 = 2885.9092
$ws.Range("M34").Value2 = -1059.1428
$ws.Range("N34").Value2 = -3289.9092
$ws.Range("H55").Value2 = 35040
$ws.Range("I55").Value2 = 20000
$ws.Range("K55").Value2 = 20000
$ws.Range("M55").Value2 = -19685
$ws.Range("H105").Value2 = 4883.4443
$ws.Range("I105").Value2 = 4883.4443
$ws.Range("K105").Value2 = 4883.4443
$ws.Range("M105").Value2 = -3136.4443
$ws.Range("H122").Value2 = 3538.926
$ws.Range("I122").Value2 = 3025.087
$ws.Range("K122").Value2 = 9075.261
$ws.Range("M122").Value2 = -6625.261
$ws.Range("H132").Value2 = 1410.8
$ws.Range("I132").Value2 = 1018
$ws.Range("K132").Value2 = 3054
$ws.Range("M132").Value2 = -524
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 158.59259
$ws.Range("I2").Value2 = 118.5
$ws.Range("J2").Value2 = 238.77777
$ws.Range("K2").Value2 = 711
$ws.Range("L2").Value2 = 1432.66662
$ws.Range("M2").Value2 = -598
$ws.Range("N2").Value2 = -1658.66662
$ws.Range("H38").Value2 = 2761.8635
$ws.Range("I38").Value2 = 291.64285
$ws.Range("J38").Value2 = 7084.75
$ws.Range("K38").Value2 = 874.9285500000001
$ws.Range("L38").Value2 = 21254.25
$ws.Range("M38").Value2 = -527.9285500000001
$ws.Range("N38").Value2 = -21948.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 15421.77
$ws.Range("I70").Value2 = 7097.6
$ws.Range("J70").Value2 = 20624.375
$ws.Range("K70").Value2 = 7097.6
$ws.Range("L70").Value2 = 20624.375
$ws.Range("M70").Value2 = -6827.6
$ws.Range("N70").Value2 = -21164.375
$ws.Range("H73").Value2 = 15421.77
$ws.Range("I73").Value2 = 7097.6
$ws.Range("J73").Value2 = 20624.375
$ws.Range("K73").Value2 = 7097.6
$ws.Range("L73").Value2 = 20624.375
$ws.Range("M73").Value2 = -6161.6
$ws.Range("N73").Value2 = -22496.375
$ws.Range("H80").Value2 = 6132.3076
$ws.Range("I80").Value2 = 3903.4285
$ws.Range("J80").Value2 = 8732.666999999999
$ws.Range("K80").Value2 = 3903.4285
$ws.Range("L80").Value2 = 8732.666999999999
$ws.Range("M80").Value2 = -2905.4285
$ws.Range("N80").Value2 = -10728.667
$ws.Range("H83").Value2 = 6132.3076
$ws.Range("I83").Value2 = 3903.4285
$ws.Range("J83").Value2 = 8732.666999999999
$ws.Range("K83").Value2 = 19517.1425
$ws.Range("L83").Value2 = 43663.335
$ws.Range("M83").Value2 = -14525.1425
$ws.Range("N83").Value2 = -53647.335
$ws.Range("H113").Value2 = 2671.5
$ws.Range("I113").Value2 = 2273
$ws.Range("K113").Value2 = 2273
$ws.Range("M113").Value2 = -103
$ws.Range("H117").Value2 = 0
$ws.Range("J117").Value2 = 0
$ws.Range("L117").Value2 = 0
$ws.Range("N117").ClearContents()
$ws.Range("H132").Value2 = 3716.1719
$ws.Range("J132").Value2 = 2361.9092
$ws.Range("L132").Value2 = 7085.7276
$ws.Range("N132").Value2 = -12145.7276
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 1946.8788
$ws.Range("I16").Value2 = 1607.0476
$ws.Range("J16").Value2 = 2541.5833
$ws.Range("K16").Value2 = 1607.0476
$ws.Range("L16").Value2 = 2541.5833
$ws.Range("M16").Value2 = -1437.0476
$ws.Range("N16").Value2 = -2881.5833
$ws.Range("H82").Value2 = 7084.6113
$ws.Range("J82").Value2 = 1716.2727
$ws.Range("L82").Value2 = 1716.2727
$ws.Range("N82").Value2 = -2438.2727
$ws.Range("H85").Value2 = 7084.6113
$ws.Range("J85").Value2 = 1716.2727
$ws.Range("L85").Value2 = 1716.2727
$ws.Range("N85").Value2 = -4212.2727
$ws.Range("H136").Value2 = 4686
$ws.Range("I136").Value2 = 4762.6943
$ws.Range("K136").Value2 = 14288.0829
$ws.Range("M136").Value2 = -11738.0829
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value2 = 12020.167
$ws.Range("I45").Value2 = 3669
$ws.Range("J45").Value2 = 12779.363
$ws.Range("K45").Value2 = 3669
$ws.Range("L45").Value2 = 12779.363
$ws.Range("M45").Value2 = -3178
$ws.Range("N45").Value2 = -13761.363
$ws.Range("H126").Value2 = 5032.8037
$ws.Range("I126").Value2 = 4505.028
$ws.Range("J126").Value2 = 6299.467
$ws.Range("K126").Value2 = 13515.084
$ws.Range("L126").Value2 = 18898.401
$ws.Range("M126").Value2 = -11045.084
$ws.Range("N126").Value2 = -23838.401
$ws.Range("H136").Value2 = 3960.9795
$ws.Range("I136").Value2 = 3890.1702
$ws.Range("K136").Value2 = 11670.5106
$ws.Range("M136").Value2 = -9120.5106
$ws.Range("H140").Value2 = 99498.5
$ws.Range("I140").Value2 = 0
$ws.Range("J140").Value2 = 99498.5
$ws.Range("K140").Value2 = 0
$ws.Range("L140").Value2 = 99498.5
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value2 = -109858.5
$ws.Range("H141").Value2 = 87000
$ws.Range("I141").Value2 = 0
$ws.Range("J141").Value2 = 87000
$ws.Range("K141").Value2 = 0
$ws.Range("L141").Value2 = 87000
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value2 = -97360
